$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" date column (C2:C5) from 45233 (2023-11-03) to
# 45243 (2023-11-13), keeping the existing date formatting intact.
$ws.Range("C2").Value = 45243
$ws.Range("C3").Value = 45243
$ws.Range("C4").Value = 45243
$ws.Range("C5").Value = 45243
